$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New volunteer-hours rows for 1-21-2018 (rows 26-29) ---
$ws.Range("A26").Value = "3:41AM 1-21-2018"
$ws.Range("B26").Value = "6:32AM 1-21-2018"
$ws.Range("C26").Value = 171

$ws.Range("A27").Value = "7:07AM 1-21-2018"
$ws.Range("B27").Value = "7:47AM 1-21-2018"
$ws.Range("B27").NumberFormat = "HH:MM:SS\ AM/PM"
$ws.Range("C27").Value = 40

$ws.Range("A28").Value = "2:07PM 1-21-2018"
$ws.Range("B28").Value = "3:04PM 1-21-2018"
$ws.Range("C28").Value = 53

$ws.Range("A29").Value = "7:08PM 1-21-2018"
$ws.Range("B29").Value = "8:00PM 1-21-2018"
$ws.Range("C29").Value = 52

# --- Move the "Total Project Hours:" summary row down to row 40 and
#     extend the SUM range to cover the newly added rows ---
$ws.Range("A40").Value = "Total Project Hours:"
$ws.Range("C40").Formula = "=SUM(C2:C39)/60"

# --- Update the view: selection ---
$ws.Range("C32").Select() | Out-Null

Write-Output "done"
